$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.303.03"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.648.32"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "520.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.73%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "2.663.63"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +9.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Value = "3.122.89"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "59.305.33"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "2.652.08"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.11%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "0.0₃0801"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.05%  "
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.896"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("E42").Value = "  +4.31%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "275.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0973"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0535"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").Value = "2.057.17"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
